$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-24 19:32:04"
$wsZh.Range("H2").Value = "2016-03-24 19:32:40"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-24 19:32:09"
$wsDe.Range("H2").Value = "2016-03-24 19:32:47"
